$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 328
    3  = 7
    4  = 10443
    6  = 948
    7  = 27
    9  = 7210
    10 = 20
    11 = 447
    13 = 130
    14 = 3214
    19 = 1040
    22 = 1656
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
